# The template used Word "complex fields" ({ m:userdoc 1/0 } and
# { m:enduserdoc }) built from w:fldChar/w:instrText runs. The parser was
# switched to TokenIteratorFieldRewriterSplit, which expects the same
# characters to appear as plain literal text runs (w:t) instead of field
# codes. This rewrites the two field paragraphs accordingly, keeping the
# _GoBack bookmark that sits in the middle of the first one.

$d = $word.ActiveDocument

function Set-ParagraphXml($paragraph, [string]$innerBodyXml) {
    $xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
           $innerBodyXml +
           '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $paragraph.Range.InsertXML($xml)
}

function Get-ParagraphIndexForPosition($doc, $pos) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $para = $doc.Paragraphs.Item($i)
        if ($pos -ge $para.Range.Start -and $pos -lt $para.Range.End) {
            return $i
        }
    }
    return 0
}

# Locate the two field-bearing paragraphs by the field's own instruction
# text (instead of a hard-coded paragraph index), using Document.Fields
# (Range.Fields is not reliably scoped to sub-ranges).
$userdocParaIndex = 0
$enduserdocParaIndex = 0
foreach ($f in $d.Fields) {
    $code = $f.Code.Text
    $idx = Get-ParagraphIndexForPosition $d $f.Code.Start
    if ($code -match "enduserdoc") {
        $enduserdocParaIndex = $idx
    } elseif ($code -match "userdoc") {
        $userdocParaIndex = $idx
    }
}

if ($userdocParaIndex -gt 0) {
    $body = '<w:p><w:r><w:t>{</w:t></w:r>' +
            '<w:r><w:t>m</w:t></w:r>' +
            '<w:r><w:t xml:space="preserve">:userdoc </w:t></w:r>' +
            '<w:r><w:t>1/0</w:t></w:r>' +
            '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
            '<w:bookmarkEnd w:id="0"/>' +
            '<w:r><w:t xml:space="preserve">}</w:t></w:r>' +
            '</w:p>'
    Set-ParagraphXml $d.Paragraphs.Item($userdocParaIndex) $body
}

if ($enduserdocParaIndex -gt 0) {
    $body = '<w:p><w:r><w:t xml:space="preserve">{m:enduserdoc}</w:t></w:r></w:p>'
    Set-ParagraphXml $d.Paragraphs.Item($enduserdocParaIndex) $body
}
